$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 14
$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 9
$ws.Range("B4").Value = 18
$ws.Range("B5").Value = 51
$ws.Range("B6").Value = 74
$ws.Range("B8").Value = 20
$ws.Range("B9").Value = 81
$ws.Range("B10").Value = 78

Write-Host "Updated cell values in column B of Sheet1."
